$wb = $excel.ActiveWorkbook


# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 481.45834
$ws.Range("I32").Value = 610.8461
$ws.Range("J32").Value = 328.54544
$ws.Range("K32").Value = 610.8461
$ws.Range("L32").Value = 328.54544
$ws.Range("M32").Value = -284.8461
$ws.Range("N32").Value = -980.54544
$ws.Range("H70").Value = 1238
$ws.Range("I70").Value = 820
$ws.Range("J70").Value = 1516.6666
$ws.Range("K70").Value = 2460
$ws.Range("L70").Value = 4549.9998
$ws.Range("M70").Value = -2190
$ws.Range("N70").Value = -5089.9998
$ws.Range("H73").Value = 1238
$ws.Range("I73").Value = 820
$ws.Range("J73").Value = 1516.6666
$ws.Range("K73").Value = 2460
$ws.Range("L73").Value = 4549.9998
$ws.Range("M73").Value = -1524
$ws.Range("N73").Value = -6421.9998
$ws.Range("H75").Value = 21978.5
$ws.Range("J75").Value = 21978.5
$ws.Range("L75").Value = 21978.5
$ws.Range("N75").Value = -23850.5
$ws.Range("H78").Value = 21978.5
$ws.Range("J78").Value = 21978.5
$ws.Range("L78").Value = 65935.5
$ws.Range("N78").Value = -75295.5
$ws.Range("H87").Value = 29363.637
$ws.Range("J87").Value = 29363.637
$ws.Range("L87").Value = 29363.637
$ws.Range("N87").Value = -31859.637
$ws.Range("H90").Value = 29363.637
$ws.Range("J90").Value = 29363.637
$ws.Range("L90").Value = 88090.91099999999
$ws.Range("N90").Value = -100570.911
$ws.Range("H96").Value = 1304.8462
$ws.Range("I96").Value = 616.75
$ws.Range("J96").Value = 2405.8
$ws.Range("K96").Value = 1850.25
$ws.Range("L96").Value = 7217.400000000001
$ws.Range("M96").Value = -477.25
$ws.Range("N96").Value = -9963.400000000001
$ws.Range("H98").Value = 791.86957
$ws.Range("I98").Value = 814.9048
$ws.Range("J98").Value = 550
$ws.Range("K98").Value = 814.9048
$ws.Range("L98").Value = 550
$ws.Range("M98").Value = 683.0952
$ws.Range("N98").Value = -3546
$ws.Range("H122").Value = 791.86957
$ws.Range("I122").Value = 814.9048
$ws.Range("J122").Value = 550
$ws.Range("K122").Value = 2444.7144
$ws.Range("L122").Value = 1650
$ws.Range("M122").Value = 5.285600000000159
$ws.Range("N122").Value = -6550
$ws.Range("H138").Value = 4839.989
$ws.Range("I138").Value = 2302.7727
$ws.Range("J138").Value = 5648.9565
$ws.Range("K138").Value = 6908.3181
$ws.Range("L138").Value = 16946.8695
$ws.Range("M138").Value = -1768.3181
$ws.Range("N138").Value = -27226.8695

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2341.4119
$ws.Range("I45").Value = 2875.5
$ws.Range("J45").Value = 1866.6666
$ws.Range("K45").Value = 2875.5
$ws.Range("L45").Value = 1866.6666
$ws.Range("M45").Value = -2498.5
$ws.Range("N45").Value = -2620.6666
$ws.Range("H63").Value = 2370.8333
$ws.Range("I63").Value = 2131.818
$ws.Range("J63").Value = 5000
$ws.Range("K63").Value = 2131.818
$ws.Range("L63").Value = 5000
$ws.Range("M63").Value = -1445.818
$ws.Range("N63").Value = -6372
$ws.Range("H66").Value = 2370.8333
$ws.Range("I66").Value = 2131.818
$ws.Range("J66").Value = 5000
$ws.Range("K66").Value = 10659.09
$ws.Range("L66").Value = 25000
$ws.Range("M66").Value = -7227.09
$ws.Range("N66").Value = -31864
$ws.Range("H74").Value = 206002.81
$ws.Range("I74").Value = 2222.6487
$ws.Range("J74").Value = 834325
$ws.Range("K74").Value = 2222.6487
$ws.Range("L74").Value = 834325
$ws.Range("M74").Value = -1348.6487
$ws.Range("N74").Value = -836073
$ws.Range("H77").Value = 206002.81
$ws.Range("I77").Value = 2222.6487
$ws.Range("J77").Value = 834325
$ws.Range("K77").Value = 11113.2435
$ws.Range("L77").Value = 4171625
$ws.Range("M77").Value = -6745.2435
$ws.Range("N77").Value = -4180361
$ws.Range("H97").Value = 1816.5588
$ws.Range("I97").Value = 1794.4333
$ws.Range("J97").Value = 1982.5
$ws.Range("K97").Value = 1794.4333
$ws.Range("L97").Value = 1982.5
$ws.Range("M97").Value = -1298.4333
$ws.Range("N97").Value = -2974.5
$ws.Range("H132").Value = 18616.842
$ws.Range("I132").Value = 22667.18
$ws.Range("J132").Value = 3038.6155
$ws.Range("K132").Value = 68001.54000000001
$ws.Range("L132").Value = 9115.8465
$ws.Range("M132").Value = -65471.54000000001
$ws.Range("N132").Value = -14175.8465

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2598.1973
$ws.Range("I31").Value = 1304.2188
$ws.Range("K31").Value = 1304.2188
$ws.Range("M31").Value = -1009.2188
$ws.Range("H34").Value = 2598.1973
$ws.Range("I34").Value = 1304.2188
$ws.Range("K34").Value = 1304.2188
$ws.Range("M34").Value = -1102.2188
$ws.Range("H99").Value = 5111.3335
$ws.Range("I99").Value = 7139.3335
$ws.Range("J99").Value = 3083.3333
$ws.Range("K99").Value = 7139.3335
$ws.Range("L99").Value = 3083.3333
$ws.Range("M99").Value = -5641.3335
$ws.Range("N99").Value = -6079.3333
$ws.Range("H126").Value = 5111.3335
$ws.Range("I126").Value = 7139.3335
$ws.Range("J126").Value = 3083.3333
$ws.Range("K126").Value = 21418.0005
$ws.Range("L126").Value = 9249.999899999999
$ws.Range("M126").Value = -18948.0005
$ws.Range("N126").Value = -14189.9999
$ws.Range("H132").Value = 2725.3125
$ws.Range("I132").Value = 2211.2856
$ws.Range("J132").Value = 3125.111
$ws.Range("K132").Value = 6633.8568
$ws.Range("L132").Value = 9375.332999999999
$ws.Range("M132").Value = -4103.8568
$ws.Range("N132").Value = -14435.333
$ws.Range("H134").Value = 4420.1904
$ws.Range("I134").Value = 4543.3687
$ws.Range("J134").Value = 3250
$ws.Range("K134").Value = 13630.1061
$ws.Range("L134").Value = 9750
$ws.Range("M134").Value = -11095.1061
$ws.Range("N134").Value = -14820

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 130.625
$ws.Range("I14").Value = 130.625
$ws.Range("K14").Value = 391.875
$ws.Range("M14").Value = -218.875
$ws.Range("H131").Value = 2021030.4
$ws.Range("J131").Value = 1274.7142
$ws.Range("L131").Value = 3824.1426
$ws.Range("N131").Value = -13904.1426
$ws.Range("H132").Value = 1054.2122
$ws.Range("I132").Value = 523.05884
$ws.Range("J132").Value = 1618.5625
$ws.Range("K132").Value = 4707.52956
$ws.Range("L132").Value = 14567.0625
$ws.Range("M132").Value = -2177.52956
$ws.Range("N132").Value = -19627.0625
$ws.Range("H137").Value = 4014.375
$ws.Range("I137").Value = 1519
$ws.Range("J137").Value = 5262.0625
$ws.Range("K137").Value = 4557
$ws.Range("L137").Value = 15786.1875
$ws.Range("M137").Value = 543
$ws.Range("N137").Value = -25986.1875

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2157.7827
$ws.Range("I40").Value = 2057.6365
$ws.Range("J40").Value = 2249.5833
$ws.Range("K40").Value = 2057.6365
$ws.Range("L40").Value = 2249.5833
$ws.Range("M40").Value = -1921.6365
$ws.Range("N40").Value = -2521.5833
$ws.Range("H92").Value = 27900
$ws.Range("J92").Value = 27900
$ws.Range("L92").Value = 27900
$ws.Range("N92").Value = -32892

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 74044.44500000001
$ws.Range("J46").Value = 74044.44500000001
$ws.Range("L46").Value = 74044.44500000001
$ws.Range("N46").Value = -74506.44500000001
$ws.Range("H58").Value = 1335
$ws.Range("I58").Value = 1335
$ws.Range("K58").Value = 1335
$ws.Range("M58").Value = -1027
$ws.Range("H63").Value = 27333.334
$ws.Range("J63").Value = 27333.334
$ws.Range("L63").Value = 27333.334
$ws.Range("N63").Value = -28581.334
$ws.Range("H66").Value = 27333.334
$ws.Range("J66").Value = 27333.334
$ws.Range("L66").Value = 82000.00199999999
$ws.Range("N66").Value = -88240.00199999999
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()
$ws.Range("H81").Value = 71430780
$ws.Range("I81").Value = 125002104
$ws.Range("J81").Value = 2363.3333
$ws.Range("K81").Value = 250004208
$ws.Range("L81").Value = 4726.6666
$ws.Range("M81").Value = -250003147
$ws.Range("N81").Value = -6848.6666
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()
$ws.Range("H84").Value = 71430780
$ws.Range("I84").Value = 125002104
$ws.Range("J84").Value = 2363.3333
$ws.Range("K84").Value = 1250021040
$ws.Range("L84").Value = 23633.333
$ws.Range("M84").Value = -1250015736
$ws.Range("N84").Value = -34241.333
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()
$ws.Range("H134").Value = 74044.44500000001
$ws.Range("J134").Value = 74044.44500000001
$ws.Range("L134").Value = 222133.335
$ws.Range("N134").Value = -227203.335
$ws.Range("H136").Value = 3184.4038
$ws.Range("I136").Value = 1030.359
$ws.Range("J136").Value = 9646.538
$ws.Range("K136").Value = 3091.077
$ws.Range("L136").Value = 28939.614
$ws.Range("M136").Value = -541.0769999999998
$ws.Range("N136").Value = -34039.614

Write-Output "Applied all profit table updates"